$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 314 (low/close changed) ---
$ws.Cells.Item(314, 5).Value = 106.43   # E314 (low)
$ws.Cells.Item(314, 6).Value = 109      # F314 (close)

# --- Append new rows 315-317 ---

# Row 315
$ws.Cells.Item(315, 1).Value = 45170.33333333334
$ws.Cells.Item(315, 2).Value = "FX_IDC:USDBDT"
$ws.Cells.Item(315, 3).Value = 109
$ws.Cells.Item(315, 4).Value = 110.28
$ws.Cells.Item(315, 5).Value = 108.53
$ws.Cells.Item(315, 6).Value = 109.97
$ws.Cells.Item(315, 7).Value = 0

# Row 316
$ws.Cells.Item(316, 1).Value = 45201.375
$ws.Cells.Item(316, 2).Value = "FX_IDC:USDBDT"
$ws.Cells.Item(316, 3).Value = 109.97
$ws.Cells.Item(316, 4).Value = 110.4
$ws.Cells.Item(316, 5).Value = 108.53
$ws.Cells.Item(316, 6).Value = 110
$ws.Cells.Item(316, 7).Value = 0

# Row 317
$ws.Cells.Item(317, 1).Value = 45231.375
$ws.Cells.Item(317, 2).Value = "FX_IDC:USDBDT"
$ws.Cells.Item(317, 3).Value = 110
$ws.Cells.Item(317, 4).Value = 110.63
$ws.Cells.Item(317, 5).Value = 108.5
$ws.Cells.Item(317, 6).Value = 110.48
$ws.Cells.Item(317, 7).Value = 0

# --- Copy the date-column formatting (style index used on A2:A314) down onto
#     the new date cells A315:A317 so they keep the same number format/border ---
$ws.Range("A314").Copy()
$ws.Range("A315:A317").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Refresh the sheet's used-range / dimension so it reflects the new rows
$ws.UsedRange | Out-Null
